$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.766.64"
$ws.Range("E2").Value = "  +3.56%  "
$ws.Range("D3").Value = "3.698.38"
$ws.Range("E3").Value = "  +7.88%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'589.33"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'180.79"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "3.689.00"
$ws.Range("E7").Value = "  +7.84%  "
$ws.Range("E8").Value = "  +4.13%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").Value = "'50.05"
$ws.Range("E12").Value = "  +3.81%  "
$ws.Range("D13").Value = "'0.0000289"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").Value = "4.292.99"
$ws.Range("E14").Value = "  +7.91%  "
$ws.Range("D15").Value = "'683.64"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'9.05"
$ws.Range("E16").Value = "  +4.60%  "
$ws.Range("D17").Value = "3.700.69"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("D18").Value = "71.758.63"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "'11.69"
$ws.Range("E21").Value = "  +3.24%  "
$ws.Range("D22").Value = "'6.49"
$ws.Range("E22").Value = "  +20.68%  "
$ws.Range("D23").Value = "'0.945"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("E24").Value = "  +4.86%  "
$ws.Range("D25").Value = "'103.91"
$ws.Range("E25").Value = "  +3.14%  "
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("E27").Value = "  +5.32%  "
$ws.Range("D28").Value = "'10.30"
$ws.Range("E28").Value = "  +6.53%  "
$ws.Range("D29").Value = "'35.60"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("D30").Value = "'9.29"
$ws.Range("E30").Value = "  +6.00%  "
$ws.Range("D31").Value = "'7.40"
$ws.Range("D32").Value = "'4.27"
$ws.Range("E32").Value = "  +14.35%  "
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("D35").Value = "'563.07"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "'59.72"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("D37").Value = "3.754.62"
$ws.Range("E37").Value = "  +3.85%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").Value = "0.0₃0784"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("D41").Value = "'35.74"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'3.47"
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("E43").Value = "  +9.93%  "
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("E45").Value = "  +5.62%  "
$ws.Range("E46").Value = "  +8.45%  "
$ws.Range("D47").Value = "'3.38"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "'136.03"
$ws.Range("E51").Value = "  +3.84%  "

Write-Host "Applied cryptos update"
